$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 2165.691048665342
$ws.Range("C2").Value = 1838.018729192084
$ws.Range("D2").Value = 1303.400515533381
$ws.Range("E2").Value = 2213.267102976351
$ws.Range("F2").Value = 2175.325749412518
$ws.Range("G2").Value = 2108.131502855095
$ws.Range("H2").Value = 2214.805403650955
$ws.Range("B3").Value = 2183.136329287854
$ws.Range("C3").Value = 1847.781862906447
$ws.Range("D3").Value = 1158.652876897095
$ws.Range("E3").Value = 2214.879550218847
$ws.Range("F3").Value = 2187.073464779258
$ws.Range("G3").Value = 2100.219135539186
$ws.Range("H3").Value = 2215.442358444167
$ws.Range("B4").Value = 2127.598171040406
$ws.Range("C4").Value = 1850.061862514113
$ws.Range("D4").Value = 1178.468385611193
$ws.Range("E4").Value = 2204.643962329801
$ws.Range("F4").Value = 2140.236464658321
$ws.Range("G4").Value = 2096.593753118568
$ws.Range("H4").Value = 2206.561388682944
$ws.Range("B5").Value = 2185.201902321746
$ws.Range("C5").Value = 1856.916359620758
$ws.Range("D5").Value = 1259.317519699002
$ws.Range("E5").Value = 2213.264139928402
$ws.Range("F5").Value = 2187.615832963864
$ws.Range("G5").Value = 2070.329385050662
$ws.Range("H5").Value = 2213.921757419303
$ws.Range("B6").Value = 2191.64712965028
$ws.Range("C6").Value = 1855.084687763432
$ws.Range("D6").Value = 440.1144253121332
$ws.Range("E6").Value = 2219.368179815974
$ws.Range("F6").Value = 2193.78184930315
$ws.Range("G6").Value = 1894.731092481184
$ws.Range("H6").Value = 2219.610231407072
$ws.Range("B7").Value = 2171.940903516458
$ws.Range("C7").Value = 1866.352806353311
$ws.Range("D7").Value = 668.115238023379
$ws.Range("E7").Value = 2209.792567099627
$ws.Range("F7").Value = 2175.768832331416
$ws.Range("G7").Value = 1916.80201782857
$ws.Range("H7").Value = 2210.892787299947
$ws.Range("B8").Value = 2172.503856724016
$ws.Range("C8").Value = 1810.082880098869
$ws.Range("D8").Value = 724.304622318827
$ws.Range("E8").Value = 2204.063040089982
$ws.Range("F8").Value = 2173.588269967609
$ws.Range("G8").Value = 1896.532347435891
$ws.Range("H8").Value = 2204.46439533886
$ws.Range("B9").Value = 2187.02243935776
$ws.Range("C9").Value = 1854.667135943324
$ws.Range("D9").Value = 702.5796034120372
$ws.Range("E9").Value = 2207.328665015756
$ws.Range("F9").Value = 2191.725020862515
$ws.Range("G9").Value = 1920.750419483241
$ws.Range("H9").Value = 2208.723841640438
$ws.Range("B10").Value = 1922.900717917014
$ws.Range("C10").Value = 1939.815918085547
$ws.Range("D10").Value = 1378.827377776219
$ws.Range("E10").Value = 2152.809230531152
$ws.Range("F10").Value = 2001.459579491783
$ws.Range("G10").Value = 2063.638446140928
$ws.Range("H10").Value = 2161.726782239577
$ws.Range("B11").Value = 1866.433802825766
$ws.Range("C11").Value = 1953.123059391344
$ws.Range("D11").Value = 1207.322159653751
$ws.Range("E11").Value = 2155.703348248459
$ws.Range("F11").Value = 1940.725141317559
$ws.Range("G11").Value = 2045.065604917479
$ws.Range("H11").Value = 2160.859079728361
$ws.Range("B12").Value = 1582.176651645972
$ws.Range("C12").Value = 1925.405581190644
$ws.Range("D12").Value = 462.4579046651045
$ws.Range("E12").Value = 2098.051720771348
$ws.Range("F12").Value = 1605.030983419213
$ws.Range("G12").Value = 1950.445220564697
$ws.Range("H12").Value = 2098.669692755373
$ws.Range("B13").Value = 1903.651813317016
$ws.Range("C13").Value = 1935.285702291976
$ws.Range("D13").Value = 1152.058979345561
$ws.Range("E13").Value = 2145.920377647311
$ws.Range("F13").Value = 1954.061397624516
$ws.Range("G13").Value = 2019.552879454643
$ws.Range("H13").Value = 2150.349025874085